# Generate Report for Handoff
# Updates the localization-status report: refreshes the "Latest HO Xliff
# Generate Date" / "Latest Handoff Datetime" timestamps for the rows that
# were just (re)handed off, and records their new "Priority" (ht).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 13)

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-24 22:21:51"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (H) + "Priority" (E) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-24 22:21:46"
}

# --- de-de sheet: "Latest Handoff Datetime" (H) + "Priority" (E) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-24 22:21:46"
}
